$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 1036, shifting the existing data (old rows
# 1036-1164) down to 1038-1166. This mirrors Excel's "insert copied cells"
# behaviour, taking formatting from the row above.
$ws.Rows("1036:1037").Insert()

# New row 1036 — Primera, fecha 45124 (2023-07-17)
$ws.Range("A1036").Value = 3
$ws.Range("B1036").Value = "Femacal de La Calera"
$ws.Range("C1036").Value = "Coquimbo"
$ws.Range("D1036").Value = 45124
$ws.Range("E1036").Value = 5
$ws.Range("F1036").Value = 100112023
$ws.Range("G1036").Value = "Brócoli"
$ws.Range("H1036").Value = "Sin especificar"
$ws.Range("I1036").Value = "Primera"
$ws.Range("J1036").Value = 3700
$ws.Range("K1036").Value = 650
$ws.Range("L1036").Value = 700
$ws.Range("M1036").Value = 676
$ws.Range("N1036").Value = "$/unidad"
$ws.Range("O1036").Value = "Provincia de Quillota"
$ws.Range("P1036").Value = 676
$ws.Range("Q1036").Value = 1
$ws.Range("R1036").Value = "Hortaliza"

# New row 1037 — Segunda, fecha 45124 (2023-07-17)
$ws.Range("A1037").Value = 3
$ws.Range("B1037").Value = "Femacal de La Calera"
$ws.Range("C1037").Value = "Coquimbo"
$ws.Range("D1037").Value = 45124
$ws.Range("E1037").Value = 5
$ws.Range("F1037").Value = 100112023
$ws.Range("G1037").Value = "Brócoli"
$ws.Range("H1037").Value = "Sin especificar"
$ws.Range("I1037").Value = "Segunda"
$ws.Range("J1037").Value = 1800
$ws.Range("K1037").Value = 500
$ws.Range("L1037").Value = 500
$ws.Range("M1037").Value = 500
$ws.Range("N1037").Value = "$/unidad"
$ws.Range("O1037").Value = "Provincia de Quillota"
$ws.Range("P1037").Value = 500
$ws.Range("Q1037").Value = 1
$ws.Range("R1037").Value = "Hortaliza"
